# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-13 09:17:10
#
# In column G ("Recorded By"), values are comma-separated lists such as
# "System, dnasr281@gmail.com" or "System, backup@backdoor.com".
# This swaps the first two tokens when the first one is "System" and the
# second one is either "dnasr281@gmail.com" or "backup@backdoor.com",
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# and "System, backup@backdoor.com, system" -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $parts = $val -split ", "

    if ($parts.Count -ge 2 -and $parts[0] -eq "System" -and ($parts[1] -eq "dnasr281@gmail.com" -or $parts[1] -eq "backup@backdoor.com")) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $newVal = [string]::Join(", ", $parts)
        $cell.Value2 = $newVal
    }
}
